# Script: applies the 01-11-2023 14:45 scraper update to the Indonesia
# Liga 1 2023-2024 betting odds sheet.
#
# The source site re-ordered a number of already-scraped fixtures (the
# "home"/"away" pair that used to sit in row N now sits in row N+1, and
# vice versa), and one brand-new fixture (Madura United vs Persib Bandung)
# was appended as the new last row (154).
#
# Columns A:E (Indice, pais, torneio, temporada, data_partida) are stable
# per physical row; only F:V (home .. url_partida) travel with the match,
# so each "reorder" below is implemented as a swap of the F:V range
# between the two affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this interpreter silently drops calls that use named parameters
# (e.g. "-RowA 22 -RowB 23") -- always call with positional arguments.
function Swap-MatchRows($RowA, $RowB) {
    $rangeA = $ws.Range("F" + $RowA + ":V" + $RowA)
    $rangeB = $ws.Range("F" + $RowB + ":V" + $RowB)

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

# Every adjacent-row pair whose fixture data was swapped.
$swapPairs = @(
    , @(22, 23)
    , @(24, 25)
    , @(28, 29)
    , @(42, 43)
    , @(44, 45)
    , @(47, 48)
    , @(51, 52)
    , @(57, 58)
    , @(60, 61)
    , @(62, 63)
    , @(74, 75)
    , @(76, 77)
    , @(82, 83)
    , @(84, 85)
    , @(96, 97)
    , @(127, 128)
    , @(136, 137)
)

foreach ($pair in $swapPairs) {
    Swap-MatchRows $pair[0] $pair[1]
}

# New fixture appended as row 154 (Indice 153). Copy formatting from the
# previous last row (153) first so styles (bold/border on A, datetime
# number format on E) match the rest of the table, then fill in values.
$ws.Range("A153:V153").Copy() | Out-Null
$ws.Range("A154:V154").PasteSpecial(-4122) | Out-Null

$ws.Range("A154").Value2 = 153
$ws.Range("B154").Value2 = "indonesia"
$ws.Range("C154").Value2 = "liga-1"
$ws.Range("D154").Value2 = "2023-2024"
$ws.Range("E154").Value2 = 45231.54166666666
$ws.Range("F154").Value2 = "Madura United"
$ws.Range("G154").Value2 = 0
$ws.Range("H154").Value2 = "Persib Bandung"
$ws.Range("I154").Value2 = 1
$ws.Range("J154").Value2 = 2.28
$ws.Range("K154").Value2 = "31/10/2023 01:12"
$ws.Range("L154").Value2 = 2.8
$ws.Range("M154").Value2 = "01/11/2023 12:59"
$ws.Range("N154").Value2 = 3.25
$ws.Range("O154").Value2 = "31/10/2023 01:12"
$ws.Range("P154").Value2 = 3.37
$ws.Range("Q154").Value2 = "01/11/2023 12:57"
$ws.Range("R154").Value2 = 2.82
$ws.Range("S154").Value2 = "31/10/2023 01:12"
$ws.Range("T154").Value2 = 2.51
$ws.Range("U154").Value2 = "01/11/2023 12:33"
$ws.Range("V154").Value2 = "https://www.betexplorer.com/football/indonesia/liga-1/madura-united-persib-bandung/GhWzUlSt/"
